# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: title / link
$ws.Range("D9").Value = "MSc AI/DS Prep F2022 시험 후기"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msc-ai-ds-prep-f2022-review/#utm_source=rss&utm_medium=rss&utm_campaign=msc-ai-ds-prep-f2022-review"

# Row 27: title / link
$ws.Range("D27").Value = "개발자를 위한 AWS 클라우드 보안 (3) - 인프라 보안과 사고 대응"
$ws.Range("E27").Value = "https://blog.pingpong.us/aws-cloud-security-for-devs-3/"

# Row 37: title / link
$ws.Range("D37").Value = "[Paper Review] Unsupervised Time-Series Representation Learning with Iterative Bilinear Temporal-Spectral Fusion"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2219&mod=document&pageid=1"

# Row 51: title / link
$ws.Range("D51").Value = "[python] 주소를 위도, 경도 좌표로 바꿔주는 Geocoder API 사용법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EC%A3%BC%EC%86%8C%EB%A5%BC-%EC%9C%84%EB%8F%84-%EA%B2%BD%EB%8F%84-%EC%A2%8C%ED%91%9C%EB%A1%9C-%EB%B0%94%EA%BF%94%EC%A3%BC%EB%8A%94-Geocoder-API-%EC%82%AC%EC%9A%A9%EB%B2%95"

# Row 52: title only
$ws.Range("D52").Value = "숨은 DS"
